$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (batsman) to make room for ownTeam and oppTeam
$ws.Range("D1:E1").EntireColumn.Insert()

# Set header values for new columns
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Set data values for new columns (row 2)
$ws.Range("D2").Value = "Sunrisers Hyderabad"
$ws.Range("E2").Value = "Kolkata Knight Riders"
